$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1:C1").ClearContents()

$ws.Range("B1").Value = "a"
$ws.Range("C1").Value = "b"
$ws.Range("D1").Value = "c"
$ws.Range("F1").Value = "d"

$ws.Range("A2").Value = "e"
$ws.Range("C2").Value = "f"
$ws.Range("D2").Value = "g"
$ws.Range("F2").Value = "h"

$ws.Range("A3").Value = "I"
$ws.Range("B3").Value = "j"
$ws.Range("D3").Value = "k"
$ws.Range("F3").Value = "l"

$ws.Range("A5").Value = "m"
$ws.Range("B5").Value = "n"
$ws.Range("F5").Value = "o"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A1:F5").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("A1").Select() | Out-Null

$ws.Activate()
$ws.Range("A1:F5").Select() | Out-Null
$ws.Range("F5").Activate()
